$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (former rows 6-9) before rewriting the remaining data rows.
$ws.Rows("6:9").Delete()

# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Sema6d"
$ws.Range("C2").Value2 = "Trem2"
$ws.Range("D2").Value2 = "Resolving-Mac"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 58.62893533333334
$ws.Range("H2").Value2 = 175.886806
$ws.Range("I2").Value2 = 0.5702456571409142
$ws.Range("J2").Value2 = 0.5702456571409142
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 83.28058766666668
$ws.Range("N2").Value2 = 249.841763
$ws.Range("O2").Value2 = 1
$ws.Range("P2").Value2 = 1
$ws.Range("Q2").Value2 = 4882.652188830998
$ws.Range("R2").Value2 = 43943.86969947898
$ws.Range("S2").Value2 = 0.5702456571409142
$ws.Range("T2").Value2 = 0.5702456571409142

# Row 3
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("B3").Value2 = "Sema6d"
$ws.Range("C3").Value2 = "Trem2"
$ws.Range("D3").Value2 = "Resolving-Mac"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 14.28901333333333
$ws.Range("H3").Value2 = 42.86704
$ws.Range("I3").Value2 = 0.1389799721218762
$ws.Range("J3").Value2 = 0.1389799721218763
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 83.28058766666668
$ws.Range("N3").Value2 = 249.841763
$ws.Range("O3").Value2 = 1
$ws.Range("P3").Value2 = 1
$ws.Range("Q3").Value2 = 1189.997427576836
$ws.Range("R3").Value2 = 10709.97684819152
$ws.Range("S3").Value2 = 0.1389799721218762
$ws.Range("T3").Value2 = 0.1389799721218763

# Row 4
$ws.Range("A4").Value2 = "MuSCs"
$ws.Range("B4").Value2 = "Sema6d"
$ws.Range("C4").Value2 = "Trem2"
$ws.Range("D4").Value2 = "Resolving-Mac"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 27.27518533333334
$ws.Range("H4").Value2 = 81.82555600000001
$ws.Range("I4").Value2 = 0.265288050953297
$ws.Range("J4").Value2 = 0.2652880509532971
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 83.28058766666668
$ws.Range("N4").Value2 = 249.841763
$ws.Range("O4").Value2 = 1
$ws.Range("P4").Value2 = 1
$ws.Range("Q4").Value2 = 2271.493463277248
$ws.Range("R4").Value2 = 20443.44116949523
$ws.Range("S4").Value2 = 0.265288050953297
$ws.Range("T4").Value2 = 0.2652880509532971

# Row 5
$ws.Range("A5").Value2 = "Resolving-Mac"
$ws.Range("B5").Value2 = "Sema6d"
$ws.Range("C5").Value2 = "Trem2"
$ws.Range("D5").Value2 = "Resolving-Mac"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 2.620337
$ws.Range("H5").Value2 = 7.861011
$ws.Range("I5").Value2 = 0.02548631978391236
$ws.Range("J5").Value2 = 0.02548631978391236
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 83.28058766666668
$ws.Range("N5").Value2 = 249.841763
$ws.Range("O5").Value2 = 1
$ws.Range("P5").Value2 = 1
$ws.Range("Q5").Value2 = 218.2232052447103
$ws.Range("R5").Value2 = 1964.008847202393
$ws.Range("S5").Value2 = 0.02548631978391236
$ws.Range("T5").Value2 = 0.02548631978391236
